$d = $word.ActiveDocument

# The first paragraph originally reads (across several differently formatted runs):
#   "Basic " + "    " + "<---" + "M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0" + "    " + "let" + " demonstration :"
# We need to remove everything between "Basic " and "let", i.e. the four spaces,
# the "<---" marker and the version-mismatch warning text, and the trailing four spaces,
# so the paragraph reads "Basic let demonstration :".

$range = $d.Content
$found = $range.Find.Execute(
    "    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0    ",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found -and $range.Text -eq "    <---M2Doc version mismatch: template is 3.1.1 and runtime is 3.2.0    ") {
    $range.Delete()
} else {
    throw "Could not locate the version-mismatch warning text to remove"
}
